$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$styleD = $ws.Range("D2").Style
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "25.851.30"
$ws.Range("D2").Style = $styleD
$ws.Range("E2").Value = "  +0.20%  "

$styleD = $ws.Range("D3").Style
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.743.03"
$ws.Range("D3").Style = $styleD
$ws.Range("E3").Value = "  -0.87%  "

$styleD = $ws.Range("D4").Style
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9999"
$ws.Range("D4").Style = $styleD
$ws.Range("E4").Value = "  -0.13%  "

$styleD = $ws.Range("D5").Style
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "225.13"
$ws.Range("D5").Style = $styleD
$ws.Range("E5").Value = "  -5.16%  "

$styleD = $ws.Range("D6").Style
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.9990"
$ws.Range("D6").Style = $styleD
$ws.Range("E6").Value = "  -0.13%  "

$styleD = $ws.Range("D7").Style
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5140"
$ws.Range("D7").Style = $styleD
$ws.Range("E7").Value = "  +1.47%  "

$styleD = $ws.Range("D8").Style
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2785"
$ws.Range("D8").Style = $styleD
$ws.Range("E8").Value = "  +4.59%  "

$styleD = $ws.Range("D9").Style
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "39.04"
$ws.Range("D9").Style = $styleD
$ws.Range("E9").Value = "  -5.62%  "

$styleD = $ws.Range("D10").Style
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.06093"
$ws.Range("D10").Style = $styleD
$ws.Range("E10").Value = "  -1.83%  "

$styleD = $ws.Range("D11").Style
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.734.98"
$ws.Range("D11").Style = $styleD
$ws.Range("E11").Value = "  -1.32%  "

$styleD = $ws.Range("D12").Style
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.06987"
$ws.Range("D12").Style = $styleD
$ws.Range("E12").Value = "  +0.65%  "

$styleD = $ws.Range("D13").Style
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "15.23"
$ws.Range("D13").Style = $styleD
$ws.Range("E13").Value = "  -2.42%  "

$styleD = $ws.Range("D14").Style
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.6351"
$ws.Range("D14").Style = $styleD
$ws.Range("E14").Value = "  +4.76%  "

$styleD = $ws.Range("D15").Style
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "4.507"
$ws.Range("D15").Style = $styleD
$ws.Range("E15").Value = "  +0.86%  "

$styleD = $ws.Range("D16").Style
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "76.50"
$ws.Range("D16").Style = $styleD
$ws.Range("E16").Value = "  -1.22%  "

$styleD = $ws.Range("D17").Style
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.9999"
$ws.Range("D17").Style = $styleD
$ws.Range("E17").Value = "  -0.11%  "

$styleD = $ws.Range("D18").Style
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "1.000"
$ws.Range("D18").Style = $styleD
$ws.Range("E18").Value = "  -0.04%  "

$styleD = $ws.Range("D19").Style
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "25.875.80"
$ws.Range("D19").Style = $styleD
$ws.Range("E19").Value = "  +0.15%  "

$styleD = $ws.Range("D20").Style
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "11.46"
$ws.Range("D20").Style = $styleD
$ws.Range("E20").Value = "  -1.51%  "

$styleD = $ws.Range("D21").Style
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.000006613"
$ws.Range("D21").Style = $styleD
$ws.Range("E21").Value = "  -3.09%  "

$styleD = $ws.Range("D22").Style
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "1.958.70"
$ws.Range("D22").Style = $styleD
$ws.Range("E22").Value = "  -0.94%  "

$styleD = $ws.Range("D23").Style
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "4.086"
$ws.Range("D23").Style = $styleD
$ws.Range("E23").Value = "  +0.44%  "

$styleD = $ws.Range("D24").Style
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "8.502"
$ws.Range("D24").Style = $styleD
$ws.Range("E24").Value = "  +4.12%  "

$styleD = $ws.Range("D25").Style
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "5.107"
$ws.Range("D25").Style = $styleD
$ws.Range("E25").Value = "  -1.73%  "

$styleD = $ws.Range("D26").Style
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "137.39"
$ws.Range("D26").Style = $styleD
$ws.Range("E26").Value = "  -0.31%  "

$styleD = $ws.Range("D27").Style
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "1.504"
$ws.Range("D27").Style = $styleD
$ws.Range("E27").Value = "  +3.15%  "

$styleD = $ws.Range("D28").Style
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.819"
$ws.Range("D28").Style = $styleD
$ws.Range("E28").Value = "  -0.11%  "

$styleD = $ws.Range("D29").Style
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "14.99"
$ws.Range("D29").Style = $styleD
$ws.Range("E29").Value = "  -0.12%  "

$styleD = $ws.Range("D30").Style
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "102.99"
$ws.Range("D30").Style = $styleD
$ws.Range("E30").Value = "  +0.30%  "

$styleD = $ws.Range("D31").Style
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.08270"
$ws.Range("D31").Style = $styleD
$ws.Range("E31").Value = "  +0.44%  "

$styleD = $ws.Range("D32").Style
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.619"
$ws.Range("D32").Style = $styleD
$ws.Range("E32").Value = "  -1.86%  "

$styleD = $ws.Range("D33").Style
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.387"
$ws.Range("D33").Style = $styleD
$ws.Range("E33").Value = "  -0.54%  "

$styleD = $ws.Range("D34").Style
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.04412"
$ws.Range("D34").Style = $styleD
$ws.Range("E34").Value = "  +0.85%  "

$styleD = $ws.Range("D35").Style
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.613"
$ws.Range("D35").Style = $styleD
$ws.Range("E35").Value = "  -1.54%  "

$styleD = $ws.Range("D36").Style
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.9698"
$ws.Range("D36").Style = $styleD
$ws.Range("E36").Value = "  -3.08%  "

$styleD = $ws.Range("D37").Style
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.5982"
$ws.Range("D37").Style = $styleD
$ws.Range("E37").Value = "  -1.51%  "

$styleD = $ws.Range("D38").Style
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.654"
$ws.Range("D38").Style = $styleD
$ws.Range("E38").Value = "  -2.65%  "

$styleD = $ws.Range("D39").Style
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01552"
$ws.Range("D39").Style = $styleD
$ws.Range("E39").Value = "  +0.43%  "

$styleD = $ws.Range("D40").Style
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.918"
$ws.Range("D40").Style = $styleD
$ws.Range("E40").Value = "  -1.00%  "

$styleD = $ws.Range("D41").Style
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.9979"
$ws.Range("D41").Style = $styleD
$ws.Range("E41").Value = "  -0.25%  "

$styleD = $ws.Range("D42").Style
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "101.01"
$ws.Range("D42").Style = $styleD
$ws.Range("E42").Value = "  -2.11%  "

$styleD = $ws.Range("D43").Style
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.3830"
$ws.Range("D43").Style = $styleD
$ws.Range("E43").Value = "  -0.10%  "

$styleD = $ws.Range("D44").Style
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.7301"
$ws.Range("D44").Style = $styleD
$ws.Range("E44").Value = "  -1.09%  "

$styleD = $ws.Range("D45").Style
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "4.863"
$ws.Range("D45").Style = $styleD
$ws.Range("E45").Value = "  -0.87%  "

$styleD = $ws.Range("D46").Style
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.05470"
$ws.Range("D46").Style = $styleD
$ws.Range("E46").Value = "  -0.39%  "

$styleD = $ws.Range("D47").Style
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "6.235"
$ws.Range("D47").Style = $styleD
$ws.Range("E47").Value = "  +4.90%  "

$styleD = $ws.Range("D48").Style
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.1104"
$ws.Range("D48").Style = $styleD
$ws.Range("E48").Value = "  +2.20%  "

$styleD = $ws.Range("D49").Style
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "52.19"
$ws.Range("D49").Style = $styleD
$ws.Range("E49").Value = "  +0.19%  "

$styleD = $ws.Range("D50").Style
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "29.68"
$ws.Range("D50").Style = $styleD
$ws.Range("E50").Value = "  -0.90%  "

$styleD = $ws.Range("D51").Style
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "7.471"
$ws.Range("D51").Style = $styleD
$ws.Range("E51").Value = "  -1.68%  "

